$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRange, $text) {
    $cellRange.NumberFormat = "@"
    $cellRange.Value = $text
    $cellRange.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "28.563.09"
Set-TextValue $ws.Range("E2") "  -2.31%  "
Set-TextValue $ws.Range("D3") "1.794.15"
Set-TextValue $ws.Range("E3") "  -1.93%  "
Set-TextValue $ws.Range("E4") "  -0.11%  "
Set-TextValue $ws.Range("D5") "231.87"
Set-TextValue $ws.Range("E5") "  -1.20%  "
Set-TextValue $ws.Range("D6") "0.5899"
Set-TextValue $ws.Range("E6") "  -1.43%  "
Set-TextValue $ws.Range("E7") "  +0.00%  "
Set-TextValue $ws.Range("D8") "0.2776"
Set-TextValue $ws.Range("E8") "  +0.69%  "
Set-TextValue $ws.Range("D9") "23.34"
Set-TextValue $ws.Range("E9") "  +0.17%  "
Set-TextValue $ws.Range("D10") "0.06763"
Set-TextValue $ws.Range("E10") "  -2.98%  "
Set-TextValue $ws.Range("D11") "0.07550"
Set-TextValue $ws.Range("E11") "  -1.08%  "
Set-TextValue $ws.Range("D12") "1.785.22"
Set-TextValue $ws.Range("E12") "  -2.43%  "
Set-TextValue $ws.Range("D13") "4.799"
Set-TextValue $ws.Range("E13") "  +0.53%  "
Set-TextValue $ws.Range("D14") "0.6139"
Set-TextValue $ws.Range("E14") "  -1.95%  "
Set-TextValue $ws.Range("D15") "2.037.08"
Set-TextValue $ws.Range("E15") "  -1.96%  "
Set-TextValue $ws.Range("D16") "75.60"
Set-TextValue $ws.Range("E16") "  -3.81%  "
Set-TextValue $ws.Range("D17") "0.000008923"
Set-TextValue $ws.Range("E17") "  -8.05%  "
Set-TextValue $ws.Range("D18") "28.555.50"
Set-TextValue $ws.Range("E18") "  -1.76%  "
Set-TextValue $ws.Range("D19") "5.420"
Set-TextValue $ws.Range("E19") "  -5.67%  "
Set-TextValue $ws.Range("E20") "  -0.06%  "
Set-TextValue $ws.Range("D21") "208.92"
Set-TextValue $ws.Range("E21") "  -5.79%  "
Set-TextValue $ws.Range("D22") "11.48"
Set-TextValue $ws.Range("E22") "  -0.82%  "
Set-TextValue $ws.Range("D23") "6.833"
Set-TextValue $ws.Range("E23") "  -1.02%  "
Set-TextValue $ws.Range("E24") "  -0.09%  "
Set-TextValue $ws.Range("D25") "152.46"
Set-TextValue $ws.Range("E25") "  -2.30%  "
Set-TextValue $ws.Range("D26") "8.087"
Set-TextValue $ws.Range("E26") "  +1.68%  "
Set-TextValue $ws.Range("E27") "  -2.44%  "
Set-TextValue $ws.Range("D28") "16.44"
Set-TextValue $ws.Range("E28") "  -0.43%  "
Set-TextValue $ws.Range("D29") "1.408"
Set-TextValue $ws.Range("E29") "  -2.66%  "
Set-TextValue $ws.Range("D30") "0.06193"
Set-TextValue $ws.Range("E30") "  -8.43%  "
Set-TextValue $ws.Range("D31") "1.422"
Set-TextValue $ws.Range("E31") "  -1.45%  "
Set-TextValue $ws.Range("D32") "3.793"
Set-TextValue $ws.Range("E32") "  -1.08%  "
Set-TextValue $ws.Range("D33") "3.767"
Set-TextValue $ws.Range("E33") "  -0.03%  "
Set-TextValue $ws.Range("D34") "1.729"
Set-TextValue $ws.Range("E34") "  +0.43%  "
Set-TextValue $ws.Range("D35") "1.046"
Set-TextValue $ws.Range("E35") "  -4.21%  "
Set-TextValue $ws.Range("D36") "0.6421"
Set-TextValue $ws.Range("D37") "2.504"
Set-TextValue $ws.Range("E37") "  -1.64%  "
Set-TextValue $ws.Range("D38") "2.703"
Set-TextValue $ws.Range("E38") "  -1.28%  "
Set-TextValue $ws.Range("D39") "0.01696"
Set-TextValue $ws.Range("E39") "  -2.43%  "
Set-TextValue $ws.Range("D40") "6.340"
Set-TextValue $ws.Range("E40") "  -2.59%  "
Set-TextValue $ws.Range("D41") "1.147.94"
Set-TextValue $ws.Range("E41") "  -3.96%  "
Set-TextValue $ws.Range("D42") "0.8738"
Set-TextValue $ws.Range("E42") "  -3.17%  "
Set-TextValue $ws.Range("E43") "  +0.05%  "
Set-TextValue $ws.Range("D44") "100.35"
Set-TextValue $ws.Range("E44") "  +0.04%  "
Set-TextValue $ws.Range("D45") "1.945.06"
Set-TextValue $ws.Range("E45") "  -1.94%  "
Set-TextValue $ws.Range("D46") "60.20"
Set-TextValue $ws.Range("E46") "  -3.34%  "
Set-TextValue $ws.Range("D47") "0.00000000109"
Set-TextValue $ws.Range("E47") "  -3.58%  "
Set-TextValue $ws.Range("D48") "1.586"
Set-TextValue $ws.Range("E48") "  +0.68%  "
Set-TextValue $ws.Range("D49") "8.367"
Set-TextValue $ws.Range("E49") "  -1.23%  "
Set-TextValue $ws.Range("D50") "0.05458"
Set-TextValue $ws.Range("E50") "  -1.02%  "
Set-TextValue $ws.Range("D51") "0.4476"
Set-TextValue $ws.Range("E51") "  -1.77%  "
